# Form the consolidated report: populate the "Absent" column (H) so that
# it reflects attendance - a student is marked Absent (1) on days where
# the "Real" (E) attendance count is 0, otherwise 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    6  = 0
    9  = 1
    13 = 0
    15 = 1
    16 = 0
    19 = 1
    20 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}
